$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values between column C (Green channel) and column D (Blue channel)
# for data rows 2 through 40, fixing the RBG -> RGB ordering.
for ($r = 2; $r -le 40; $r++) {
    $cCell = $ws.Range("C$r")
    $dCell = $ws.Range("D$r")
    $cVal = $cCell.Value2
    $dVal = $dCell.Value2
    $cCell.Value2 = $dVal
    $dCell.Value2 = $cVal
}
